$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "URL" column (J) ---
$ws.Columns.Item(10).ColumnWidth = 64.83203125

# Header cells (J1:J2) use the same centered style as the other header cells
$ws.Range("J1:J2").HorizontalAlignment = -4108
$ws.Range("J1").Value = "URL"
$ws.Range("J1:J2").Merge() | Out-Null

# --- Row 3: LED ---
$ws.Range("D3").Value = "COM-09661"
$ws.Range("J3").Value = "https://www.sparkfun.com/products/9661"

# --- Row 4: 10K Trimpot ---
$ws.Range("J4").Value = "https://www.sparkfun.com/products/9806"

# --- Row 5: Microcontroller ---
$ws.Range("J5").Value = "https://www.mouser.com/ProductDetail/Microchip-Technology-Atmel/ATMEGA328P-PU?qs=sGAEpiMZZMtVoztFdqDXO6rEZqxeooRg"

# --- Row 6: 16MHz Ceramic Resonator ---
$ws.Range("J6").Value = "https://www.mouser.com/ProductDetail/ABRACON/AWCR-1600MD?qs=%2fha2pyFadujFfudKd%2fAEbE32MTGIgZdnwSVrtFEweNr%2f82BKUq3Fzw%3d%3d"

# --- Row 7: 0.1uF Ceramic Capacitor ---
$ws.Range("J7").Value = "https://www.mouser.com/ProductDetail/Murata-Electronics/RDER71H104K0K1H03B?qs=%2fha2pyFadugsNiSzM4QtoPEqfGceYo24BOngKj3vXcmXztcZe0j46Uc%252bvWBqE1aS"

# --- Row 8: 10 Ohm resistor network (part number + price update) ---
$ws.Range("F8").Value = "652-4607X-1LF-330"
$ws.Range("C8").Value = "4607X-101-331LF"
$ws.Range("G8").Value = 0.255
$ws.Range("J8").Value = "https://www.mouser.com/ProductDetail/Bourns/4607X-101-331LF?qs=sGAEpiMZZMvrmc6UYKmaNWhNOGonlUnMh5dMy1XYfHQ%3d"

# --- Row 9: Tactile Switch ---
$ws.Range("J9").Value = "https://www.mouser.com/ProductDetail/Panasonic/EVQ-PE604T?qs=%2fha2pyFadui45bz44%252bGA9GlcgqdKKtsWW4AUmDdzb54ngS0CtjvNBw%3d%3d"

# --- Row 10: 10K Resistor ---
$ws.Range("J10").Value = "https://www.mouser.com/ProductDetail/Yageo/MFR-25FBF52-10K?qs=sGAEpiMZZMu61qfTUdNhG0IXHLFuiNnd4ZfMuxLN9bg%3d"

# --- Row 11: 1x6 Header Pins ---
$ws.Range("J11").Value = "https://www.mouser.com/ProductDetail/Molex/22-28-4060?qs=%2fha2pyFaduje7iG0C5h0B%2fdz9lYrXiF%2fSqnA9mKTPRg%2f8JanJs%2fc5A%3d%3d"

# --- Row 12: 2x3 Header Pins ---
$ws.Range("J12").Value = "https://www.mouser.com/ProductDetail/Molex/10-89-7061?qs=%2fha2pyFaduhUTuKyE5ihpKT%252bKedbNfhSu6VwpKmAGXgZ83GtknV09A%3d%3d"

# --- Selection matches final cursor position ---
$ws.Range("J12").Select() | Out-Null
